$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $text) {
    $r = $ws.Range($addr)
    $r.Formula = "'" + $text
    $r.Style = "Normal"
}

# Row 2 (Bitcoin)
Set-TextCell "D2" "29.127.48"
$ws.Range("E2").Value = "  -2.02%  "

# Row 3 (Ethereum)
Set-TextCell "D3" "1.852.35"
$ws.Range("E3").Value = "  -0.91%  "

# Row 4 (TetherUSD)
$ws.Range("E4").Value = "  +0.16%  "

# Row 5 (BNB)
Set-TextCell "D5" "238.97"
$ws.Range("E5").Value = "  -0.86%  "

# Row 6 (XRP)
Set-TextCell "D6" "0.6937"
$ws.Range("E6").Value = "  -5.11%  "

# Row 7 (USDC)
$ws.Range("E7").Value = "  +0.13%  "

# Row 8 (Dogecoin)
Set-TextCell "D8" "0.07754"
$ws.Range("E8").Value = "  +9.29%  "

# Row 9 (Cardano)
Set-TextCell "D9" "0.3040"
$ws.Range("E9").Value = "  -2.90%  "

# Row 10 (Solana)
Set-TextCell "D10" "23.42"
$ws.Range("E10").Value = "  -4.03%  "

# Row 11 (TRON)
Set-TextCell "D11" "0.08130"
$ws.Range("E11").Value = "  -1.28%  "

# Row 12: was WrappedEther -> now Polygon
$ws.Range("B12").Value = "Polygon"
$ws.Range("C12").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextCell "D12" "0.7272"
$ws.Range("E12").Value = "  -2.53%  "

# Row 13: was Polygon -> now WrappedEther
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextCell "D13" "1.842.40"
$ws.Range("E13").Value = "  -1.46%  "

# Row 14 (Polkadot)
Set-TextCell "D14" "5.235"
$ws.Range("E14").Value = "  -1.60%  "

# Row 15 (Litecoin)
Set-TextCell "D15" "89.14"
$ws.Range("E15").Value = "  -3.44%  "

# Row 16 (WrappedBTC)
Set-TextCell "D16" "29.131.78"
$ws.Range("E16").Value = "  -1.99%  "

# Row 17 (Uniswap)
Set-TextCell "D17" "5.794"
$ws.Range("E17").Value = "  -3.85%  "

# Row 18 (ShibaInu)
Set-TextCell "D18" "0.000007803"
$ws.Range("E18").Value = "  -0.05%  "

# Row 19 (Avalanche)
$ws.Range("E19").Value = "  -1.25%  "

# Row 20 (BitcoinCash)
Set-TextCell "D20" "237.06"
$ws.Range("E20").Value = "  -4.55%  "

# Row 21 (Dai)
Set-TextCell "D21" "0.9995"

# Row 22 (WrappedliquidstakedEther2.0)
Set-TextCell "D22" "2.096.53"
$ws.Range("E22").Value = "  -0.56%  "

# Row 23 (BinanceUSD)
$ws.Range("E23").Value = "  +0.18%  "

# Row 24 (Chainlink)
Set-TextCell "D24" "7.610"
$ws.Range("E24").Value = "  -1.17%  "

# Row 25 (Cosmos)
Set-TextCell "D25" "9.012"
$ws.Range("E25").Value = "  -1.70%  "

# Row 26 (Monero)
Set-TextCell "D26" "161.48"
$ws.Range("E26").Value = "  -0.90%  "

# Row 27 (Stellar)
Set-TextCell "D27" "0.1447"
$ws.Range("E27").Value = "  -5.83%  "

# Row 28 (EthereumClassic)
Set-TextCell "D28" "18.09"
$ws.Range("E28").Value = "  -2.38%  "

# Row 29 (LidoDAOToken)
Set-TextCell "D29" "1.981"
$ws.Range("E29").Value = "  -1.91%  "

# Row 30 (Toncoin)
$ws.Range("E30").Value = "  -2.48%  "

# Row 31 (Filecoin)
Set-TextCell "D31" "4.480"
$ws.Range("E31").Value = "  -1.37%  "

# Row 32 (PancakeSwap)
Set-TextCell "D32" "1.494"
$ws.Range("E32").Value = "  -2.12%  "

# Row 33 (InternetComputer(DFINITY))
$ws.Range("E33").Value = "  -4.33%  "

# Row 34 (Hedera)
Set-TextCell "D34" "0.05235"
$ws.Range("E34").Value = "  -0.90%  "

# Row 35 (ARBITRUM)
Set-TextCell "D35" "1.190"
$ws.Range("E35").Value = "  -3.58%  "

# Row 36 (Frax)
Set-TextCell "D36" "1.032"
$ws.Range("E36").Value = "  +3.37%  "

# Row 37 (ImmutableX)
Set-TextCell "D37" "0.7035"
$ws.Range("E37").Value = "  -6.67%  "

# Row 38 (HuobiToken)
Set-TextCell "D38" "2.654"
$ws.Range("E38").Value = "  -1.69%  "

# Row 39 (VeChain)
Set-TextCell "D39" "0.01859"
$ws.Range("E39").Value = "  -3.73%  "

# Row 40 (MXToken)
$ws.Range("E40").Value = "  -2.20%  "

# Row 41 (TrustWalletToken)
Set-TextCell "D41" "0.9333"
$ws.Range("E41").Value = "  +7.68%  "

# Row 42 (FraxShare)
Set-TextCell "D42" "6.041"
$ws.Range("E42").Value = "  +0.79%  "

# Row 43 (Maker)
Set-TextCell "D43" "1.079.74"
$ws.Range("E43").Value = "  +1.56%  "

# Row 44 (TheSandbox)
Set-TextCell "D44" "0.4275"
$ws.Range("E44").Value = "  -4.47%  "

# Row 45 (Aave)
Set-TextCell "D45" "70.52"
$ws.Range("E45").Value = "  -1.05%  "

# Row 46 (PaxDollar)
$ws.Range("E46").Value = "  +0.01%  "

# Row 47 (Quant)
Set-TextCell "D47" "102.95"
$ws.Range("E47").Value = "  -0.96%  "

# Row 48 (RenderToken)
Set-TextCell "D48" "1.783"
$ws.Range("E48").Value = "  -2.23%  "

# Row 49 (RocketPoolETH)
Set-TextCell "D49" "1.991.60"
$ws.Range("E49").Value = "  -1.21%  "

# Row 50 (EnergySwap)
Set-TextCell "D50" "9.217"
$ws.Range("E50").Value = "  -2.95%  "

# Row 51 (Aptos)
Set-TextCell "D51" "7.021"
$ws.Range("E51").Value = "  -6.15%  "
